$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197 (shifts existing rows 197-209 down to 198-210)
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new weekly price record
$ws.Cells.Item(197, 1).Value  = 3
$ws.Cells.Item(197, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(197, 3).Value  = "Coquimbo"
$ws.Cells.Item(197, 4).Value  = 44516
$ws.Cells.Item(197, 5).Value  = 5
$ws.Cells.Item(197, 6).Value  = 100112039
$ws.Cells.Item(197, 7).Value  = "Ciboulette"
$ws.Cells.Item(197, 8).Value  = "Sin especificar"
$ws.Cells.Item(197, 9).Value  = "Primera"
$ws.Cells.Item(197, 10).Value = 150
$ws.Cells.Item(197, 11).Value = 1500
$ws.Cells.Item(197, 12).Value = 1500
$ws.Cells.Item(197, 13).Value = 1500
$ws.Cells.Item(197, 14).Value = "`$/docena de atados"
$ws.Cells.Item(197, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(197, 16).Value = 500
$ws.Cells.Item(197, 17).Value = 3
$ws.Cells.Item(197, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format as the other Fecha cells
$ws.Cells.Item(197, 4).NumberFormat = $ws.Cells.Item(198, 4).NumberFormat
